$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")
$ws.Activate()

# Update the AUM figure for "3 Banken-Generali Investment-Gesellschaft mbH" (row 5)
$ws.Range("B5").Value = 378800

# Update the percentage for "Schroder Investment Management (Switzerland) AG" (row 9)
$ws.Range("C9").Value = 0.0019

# Leave the cursor where the author left it
$ws.Range("F12").Select()
